$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "2025-07-14 Monday" "2025-07-15 Tuesday"

Replace-Text "36×66=2376" "85×90=7650"
Replace-Text "23×84=1932" "69×85=5865"
Replace-Text "56×83=4648" "96×18=1728"
Replace-Text "60×53=3180" "90×21=1890"
Replace-Text "40×43=1720" "71×66=4686"

Replace-Text "42×77=3234" "21×81=1701"
Replace-Text "26×55=1430" "12×75=900"
Replace-Text "51×86=4386" "16×55=880"
Replace-Text "22×86=1892" "56×64=3584"
Replace-Text "19×54=1026" "43×82=3526"

Replace-Text "25×34=850" "20×82=1640"
Replace-Text "53×23=1219" "36×40=1440"
Replace-Text "35×71=2485" "29×81=2349"
Replace-Text "49×57=2793" "62×95=5890"
Replace-Text "32×49=1568" "13×64=832"

Replace-Text "22×21=462" "59×93=5487"
Replace-Text "53×51=2703" "53×98=5194"
Replace-Text "41×67=2747" "54×65=3510"
Replace-Text "60×31=1860" "14×81=1134"
Replace-Text "44×58=2552" "85×97=8245"

Replace-Text "17×63=1071" "67×79=5293"
Replace-Text "17×76=1292" "28×87=2436"
Replace-Text "90×84=7560" "16×92=1472"
Replace-Text "79×99=7821" "28×77=2156"
Replace-Text "15×37=555" "90×98=8820"

Write-Output "Done"
